# "Unified Database Folder in Dropbox"
# The DropboxFolder row (row 5) pointed several machines at the old
# ...\LivemRNAData folder. Point them at the new shared ...\LivemRNADatabase
# folder instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "C:\Users\albertl\Dropbox\LivemRNADatabase"
$ws.Range("C5").Value = "C:\Users\hgarcia\Documents\Dropbox\LivemRNADatabase"
$ws.Range("D5").Value = "C:\Users\Albert\Dropbox\LivemRNADatabase"
$ws.Range("E5").Value = "C:\users\hgarcia\Documents\Dropbox\LivemRNADatabase"
$ws.Range("F5").Value = "C:\Users\hgarcia\Dropbox\LivemRNADatabase"
$ws.Range("G5").Value = "C:\Users\hgarcia\Dropbox\LivemRNADatabase"

# Restore/approximate the split-pane view state: the sheet is split
# vertically after column E (so column F starts the right-hand pane),
# with the left pane focused on B5 and the right pane focused on G5.
$win = $excel.ActiveWindow
$win.SplitColumn = 5
$win.SplitRow = 0

$ws.Range("B5").Select()
$ws.Range("G5").Select()
